$d = $word.ActiveDocument

# --- Edit 1 ---------------------------------------------------------------
# The <div> block near the top of the document spells "<id>p156v_1</id>"
# across three separate runs ("<id>", "p156v_1", "</id>"). Collapse them
# into a single run holding the whole string. Word's Find & Replace
# coalesces the matched span into one run, carrying over the formatting
# of the run the match starts in (Courier New / color 7f6000 / sz 18),
# which is exactly what we want here.
$d.Content.Find.Execute("<id>p156v_1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p156v_1</id>", 2)

# --- Edit 2 -----------------------------------------------------------------
# Drop the stray trailing space that sits after "<figure>" in the figure
# block further down, right before the "<id>fig_p156v_1</id>" paragraph.
$d.Content.Find.Execute("<figure> ", $true, $false, $false, $false, $false, $true, 1, $false, "<figure>", 2)

# --- Edit 3 & 4 ---------------------------------------------------------------
# Give the "<id>fig_p156v_1</id>" paragraph and the following
# "<link>...</link>" paragraph explicit KeepWithNext/KeepTogether = False
# formatting (<w:keepNext w:val="0"/> / <w:keepLines w:val="0"/>).
# (Paragraph.Range.Text includes the trailing paragraph mark, so trim it
# before comparing.)
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd("`r")
    if ($t -eq "<id>fig_p156v_1</id>" -or $t -eq "<link>https://drive.google.com/open?id=0B9-oNrvWdlO5VjFPeWlJc05CbDQ</link>") {
        $p.KeepWithNext = $false
        $p.KeepTogether = $false
    }
}
